$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.325956666666666
$ws.Range("H2").Value = 18.97787
$ws.Range("I2").Value = 0.4468357575736242
$ws.Range("J2").Value = 0.4592138460625664
$ws.Range("M2").Value = 30.46625333333334
$ws.Range("N2").Value = 91.39876000000001
$ws.Range("O2").Value = 0.2185380492512374
$ws.Range("P2").Value = 0.2331534018544084
$ws.Range("Q2").Value = 192.7281983823556
$ws.Range("R2").Value = 1734.5537854412
$ws.Range("S2").Value = 0.09765061479583867
$ws.Range("T2").Value = 0.107067270388134
$ws.Range("G3").Value = 6.325956666666666
$ws.Range("H3").Value = 18.97787
$ws.Range("I3").Value = 0.4468357575736242
$ws.Range("J3").Value = 0.4592138460625664
$ws.Range("O3").Value = 0.2491807703757967
$ws.Range("P3").Value = 0.2658454419670822
$ws.Range("Q3").Value = 219.7519430167733
$ws.Range("R3").Value = 1977.76748715096
$ws.Range("S3").Value = 0.1113428783036484
$ws.Range("T3").Value = 0.1220799078639066
$ws.Range("G4").Value = 6.325956666666666
$ws.Range("H4").Value = 18.97787
$ws.Range("I4").Value = 0.4468357575736242
$ws.Range("J4").Value = 0.4592138460625664
$ws.Range("M4").Value = 23.69037333333334
$ws.Range("N4").Value = 71.07112000000001
$ws.Range("O4").Value = 0.1699338582153697
$ws.Range("P4").Value = 0.181298667526812
$ws.Range("Q4").Value = 149.8642751238222
$ws.Range("R4").Value = 1348.7784761144
$ws.Range("S4").Value = 0.07593252427307358
$ws.Range("T4").Value = 0.08325485840100585
$ws.Range("G5").Value = 6.325956666666666
$ws.Range("H5").Value = 18.97787
$ws.Range("I5").Value = 0.4468357575736242
$ws.Range("J5").Value = 0.4592138460625664
$ws.Range("M5").Value = 26.2168665
$ws.Range("N5").Value = 52.433733
$ws.Range("O5").Value = 0.18805669340777
$ws.Range("P5").Value = 0.1337556791894743
$ws.Range("Q5").Value = 165.846761414785
$ws.Range("R5").Value = 995.08056848871
$ws.Range("S5").Value = 0.08403045506565168
$ws.Range("T5").Value = 0.06142245987330927
$ws.Range("G6").Value = 6.325956666666666
$ws.Range("H6").Value = 18.97787
$ws.Range("I6").Value = 0.4468357575736242
$ws.Range("J6").Value = 0.4592138460625664
$ws.Range("M6").Value = 24.297748
$ws.Range("N6").Value = 72.893244
$ws.Range("O6").Value = 0.1742906287498262
$ws.Range("P6").Value = 0.1859468094622229
$ws.Range("Q6").Value = 153.7065009455866
$ws.Range("R6").Value = 1383.35850851028
$ws.Range("S6").Value = 0.07787928513541188
$ws.Range("T6").Value = 0.08538934953621061
$ws.Range("G7").Value = 6.686451000000001
$ws.Range("I7").Value = 0.4722993778644153
$ws.Range("J7").Value = 0.4853828506917099
$ws.Range("M7").Value = 30.46625333333334
$ws.Range("N7").Value = 91.39876000000001
$ws.Range("O7").Value = 0.2185380492512374
$ws.Range("P7").Value = 0.2331534018544084
$ws.Range("Q7").Value = 203.71111006692
$ws.Range("R7").Value = 1833.39999060228
$ws.Range("S7").Value = 0.1032153847010624
$ws.Range("T7").Value = 0.1131686628405626
$ws.Range("G8").Value = 6.686451000000001
$ws.Range("I8").Value = 0.4722993778644153
$ws.Range("J8").Value = 0.4853828506917099
$ws.Range("O8").Value = 0.2491807703757967
$ws.Range("P8").Value = 0.2658454419670822
$ws.Range("S8").Value = 0.1176879228242645
$ws.Range("T8").Value = 0.1290368184653799
$ws.Range("G9").Value = 6.686451000000001
$ws.Range("I9").Value = 0.4722993778644153
$ws.Range("J9").Value = 0.4853828506917099
$ws.Range("M9").Value = 23.69037333333334
$ws.Range("N9").Value = 71.07112000000001
$ws.Range("O9").Value = 0.1699338582153697
$ws.Range("P9").Value = 0.181298667526812
$ws.Range("Q9").Value = 158.40452046504
$ws.Range("R9").Value = 1425.64068418536
$ws.Range("S9").Value = 0.0802596555132189
$ws.Range("T9").Value = 0.08799926407077255
$ws.Range("G10").Value = 6.686451000000001
$ws.Range("I10").Value = 0.4722993778644153
$ws.Range("J10").Value = 0.4853828506917099
$ws.Range("M10").Value = 26.2168665
$ws.Range("N10").Value = 52.433733
$ws.Range("O10").Value = 0.18805669340777
$ws.Range("P10").Value = 0.1337556791894743
$ws.Range("Q10").Value = 175.2977932257915
$ws.Range("R10").Value = 1051.786759354749
$ws.Range("S10").Value = 0.08881905929972886
$ws.Range("T10").Value = 0.06492271286119286
$ws.Range("G11").Value = 6.686451000000001
$ws.Range("I11").Value = 0.4722993778644153
$ws.Range("J11").Value = 0.4853828506917099
$ws.Range("M11").Value = 24.297748
$ws.Range("N11").Value = 72.893244
$ws.Range("O11").Value = 0.1742906287498262
$ws.Range("P11").Value = 0.1859468094622229
$ws.Range("Q11").Value = 162.465701412348
$ws.Range("R11").Value = 1462.191312711132
$ws.Range("S11").Value = 0.0823173555261407
$ws.Range("T11").Value = 0.09025539245380199
$ws.Range("G12").Value = 1.1448225
$ws.Range("H12").Value = 2.289645
$ws.Range("I12").Value = 0.08086486456196039
$ws.Range("J12").Value = 0.05540330324572383
$ws.Range("M12").Value = 30.46625333333334
$ws.Range("N12").Value = 91.39876000000001
$ws.Range("O12").Value = 0.2185380492512374
$ws.Range("P12").Value = 0.2331534018544084
$ws.Range("Q12").Value = 34.87845230670001
$ws.Range("R12").Value = 209.2707138402
$ws.Range("S12").Value = 0.01767204975433635
$ws.Range("T12").Value = 0.0129174686257119
$ws.Range("G13").Value = 1.1448225
$ws.Range("H13").Value = 2.289645
$ws.Range("I13").Value = 0.08086486456196039
$ws.Range("J13").Value = 0.05540330324572383
$ws.Range("O13").Value = 0.2491807703757967
$ws.Range("P13").Value = 0.2658454419670822
$ws.Range("Q13").Value = 39.76899970086
$ws.Range("R13").Value = 238.61399820516
$ws.Range("S13").Value = 0.02014996924788375
$ws.Range("T13").Value = 0.01472871563779573
$ws.Range("G14").Value = 1.1448225
$ws.Range("H14").Value = 2.289645
$ws.Range("I14").Value = 0.08086486456196039
$ws.Range("J14").Value = 0.05540330324572383
$ws.Range("M14").Value = 23.69037333333334
$ws.Range("N14").Value = 71.07112000000001
$ws.Range("O14").Value = 0.1699338582153697
$ws.Range("P14").Value = 0.181298667526812
$ws.Range("Q14").Value = 27.1212724254
$ws.Range("R14").Value = 162.7276345524
$ws.Range("S14").Value = 0.01374167842907725
$ws.Range("T14").Value = 0.01004454505503363
$ws.Range("G15").Value = 1.1448225
$ws.Range("H15").Value = 2.289645
$ws.Range("I15").Value = 0.08086486456196039
$ws.Range("J15").Value = 0.05540330324572383
$ws.Range("M15").Value = 26.2168665
$ws.Range("N15").Value = 52.433733
$ws.Range("O15").Value = 0.18805669340777
$ws.Range("P15").Value = 0.1337556791894743
$ws.Range("Q15").Value = 30.01365864869625
$ws.Range("R15").Value = 120.054634594785
$ws.Range("S15").Value = 0.01520717904238943
$ws.Range("T15").Value = 0.007410506454972198
$ws.Range("G16").Value = 1.1448225
$ws.Range("H16").Value = 2.289645
$ws.Range("I16").Value = 0.08086486456196039
$ws.Range("J16").Value = 0.05540330324572383
$ws.Range("M16").Value = 24.297748
$ws.Range("N16").Value = 72.893244
$ws.Range("O16").Value = 0.1742906287498262
$ws.Range("P16").Value = 0.1859468094622229
$ws.Range("Q16").Value = 27.81660860973
$ws.Range("R16").Value = 166.89965165838
$ws.Range("S16").Value = 0.01409398808827362
$ws.Range("T16").Value = 0.01030206747221037
